$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 943.6389
$ws.Range("I17").Value = 656.3571
$ws.Range("J17").Value = 1126.4546
$ws.Range("K17").Value = 1969.0713
$ws.Range("L17").Value = 3379.3638
$ws.Range("M17").Value = -1801.0713
$ws.Range("N17").Value = -3715.3638
$ws.Range("H32").Value = 5875
$ws.Range("I32").Value = 4000
$ws.Range("J32").Value = 6500
$ws.Range("K32").Value = 4000
$ws.Range("L32").Value = 6500
$ws.Range("M32").Value = -3674
$ws.Range("N32").Value = -7152
$ws.Range("H51").Value = 21599.584
$ws.Range("I51").Value = 17833.334
$ws.Range("J51").Value = 22855
$ws.Range("K51").Value = 17833.334
$ws.Range("L51").Value = 22855
$ws.Range("M51").Value = -17349.334
$ws.Range("N51").Value = -23823
$ws.Range("H100").Value = 1982.5625
$ws.Range("I100").Value = 1973.9
$ws.Range("J100").Value = 1997
$ws.Range("K100").Value = 1973.9
$ws.Range("L100").Value = 1997
$ws.Range("M100").Value = -1432.9
$ws.Range("N100").Value = -3079
$ws.Range("H132").Value = 2405
$ws.Range("I132").Value = 952.125
$ws.Range("J132").Value = 9378.8
$ws.Range("K132").Value = 2856.375
$ws.Range("L132").Value = 28136.4
$ws.Range("M132").Value = -326.375
$ws.Range("N132").Value = -33196.39999999999
$ws.Range("H137").Value = 5335.643
$ws.Range("I137").Value = 7285.857
$ws.Range("K137").Value = 21857.571
$ws.Range("M137").Value = -19307.571
$ws.Range("H138").Value = 3462.0322
$ws.Range("I138").Value = 2635.6365
$ws.Range("J138").Value = 5482.1113
$ws.Range("K138").Value = 7906.9095
$ws.Range("L138").Value = 16446.3339
$ws.Range("M138").Value = -2766.9095
$ws.Range("N138").Value = -26726.3339

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 171180.42
$ws.Range("I32").Value = 789.7843
$ws.Range("K32").Value = 789.7843
$ws.Range("M32").Value = -502.7843
$ws.Range("H61").Value = 6671.8184
$ws.Range("I61").Value = 6839
$ws.Range("K61").Value = 6839
$ws.Range("M61").Value = -6627
$ws.Range("H102").Value = 13898451
$ws.Range("I102").Value = 18522268
$ws.Range("K102").Value = 18522268
$ws.Range("M102").Value = -18520646
$ws.Range("H132").Value = 1500
$ws.Range("I132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("M132").ClearContents()
$ws.Range("H136").Value = 6671.8184
$ws.Range("I136").Value = 6839
$ws.Range("K136").Value = 20517
$ws.Range("M136").Value = -17967
$ws.Range("H139").Value = 40000
$ws.Range("J139").Value = 40000
$ws.Range("L139").Value = 40000
$ws.Range("N139").Value = -50280

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 942.6
$ws.Range("I94").Value = 998.5
$ws.Range("K94").Value = 998.5
$ws.Range("M94").Value = -547.5
$ws.Range("H107").Value = 1658.3636
$ws.Range("I107").Value = 1741.125
$ws.Range("K107").Value = 1741.125
$ws.Range("M107").Value = 178.875

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 18520374
$ws.Range("I16").Value = 21740442
$ws.Range("J16").Value = 4984.5
$ws.Range("K16").Value = 21740442
$ws.Range("L16").Value = 4984.5
$ws.Range("M16").Value = -21740155
$ws.Range("N16").Value = -5558.5
$ws.Range("H31").Value = 5003.3335
$ws.Range("I31").Value = 4801.6
$ws.Range("K31").Value = 4801.6
$ws.Range("M31").Value = -4506.6
$ws.Range("H34").Value = 5003.3335
$ws.Range("I34").Value = 4801.6
$ws.Range("K34").Value = 4801.6
$ws.Range("M34").Value = -4599.6
$ws.Range("H107").Value = 15626102
$ws.Range("I107").Value = 50000760
$ws.Range("J107").Value = 1257.0454
$ws.Range("K107").Value = 50000760
$ws.Range("L107").Value = 1257.0454
$ws.Range("M107").Value = -49998840
$ws.Range("N107").Value = -5097.0454
$ws.Range("H113").Value = 18520374
$ws.Range("I113").Value = 21740442
$ws.Range("J113").Value = 4984.5
$ws.Range("K113").Value = 21740442
$ws.Range("L113").Value = 4984.5
$ws.Range("M113").Value = -21738272
$ws.Range("N113").Value = -9324.5
$ws.Range("H132").Value = 4109.143
$ws.Range("I132").Value = 3127.6667
$ws.Range("K132").Value = 9383.000100000001
$ws.Range("M132").Value = -6853.000100000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 84214.336
$ws.Range("I33").Value = 961.1818
$ws.Range("K33").Value = 5767.0908
$ws.Range("M33").Value = -5484.0908
$ws.Range("H75").Value = 1264
$ws.Range("I75").Value = 806.5
$ws.Range("J75").Value = 1416.5
$ws.Range("K75").Value = 2419.5
$ws.Range("L75").Value = 4249.5
$ws.Range("M75").Value = -1421.5
$ws.Range("N75").Value = -6245.5
$ws.Range("H78").Value = 1264
$ws.Range("I78").Value = 806.5
$ws.Range("J78").Value = 1416.5
$ws.Range("K78").Value = 7258.5
$ws.Range("L78").Value = 12748.5
$ws.Range("M78").Value = -2266.5
$ws.Range("N78").Value = -22732.5
$ws.Range("H107").Value = 200699.8
$ws.Range("J107").Value = 200699.8
$ws.Range("L107").Value = 602099.3999999999
$ws.Range("N107").Value = -605939.3999999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H10").Value = 25000
$ws.Range("I10").Value = 0
$ws.Range("J10").Value = 25000
$ws.Range("K10").Value = 0
$ws.Range("L10").Value = 25000
$ws.Range("M10").ClearContents()
$ws.Range("N10").Value = -25338
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 0
$ws.Range("K20").Value = 0
$ws.Range("M20").ClearContents()
$ws.Range("H43").Value = 14999.167
$ws.Range("I43").Value = 10000
$ws.Range("K43").Value = 10000
$ws.Range("M43").Value = -9849
$ws.Range("H126").Value = 4877.75
$ws.Range("J126").Value = 5000
$ws.Range("L126").Value = 15000
$ws.Range("N126").Value = -19940
$ws.Range("H132").Value = 2140.6667
$ws.Range("I132").Value = 1668.8
$ws.Range("K132").Value = 5006.4
$ws.Range("M132").Value = -2476.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 7225.9
$ws.Range("I16").Value = 12821.8
$ws.Range("K16").Value = 12821.8
$ws.Range("M16").Value = -12651.8
$ws.Range("H22").Value = 899.8
$ws.Range("I22").Value = 833
$ws.Range("J22").Value = 1000
$ws.Range("K22").Value = 833
$ws.Range("L22").Value = 1000
$ws.Range("M22").Value = -538
$ws.Range("N22").Value = -1590
$ws.Range("H27").Value = 899.8
$ws.Range("I27").Value = 833
$ws.Range("J27").Value = 1000
$ws.Range("K27").Value = 833
$ws.Range("L27").Value = 1000
$ws.Range("M27").Value = -726
$ws.Range("N27").Value = -1214
$ws.Range("H100").Value = 14499
$ws.Range("J100").Value = 14499
$ws.Range("L100").Value = 14499
$ws.Range("N100").Value = -15581
$ws.Range("H136").Value = 5665
$ws.Range("I136").Value = 5188.5713
$ws.Range("J136").Value = 6498.75
$ws.Range("K136").Value = 15565.7139
$ws.Range("L136").Value = 19496.25
$ws.Range("M136").Value = -13015.7139
$ws.Range("N136").Value = -24596.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 1433600
$ws.Range("I2").Value = 2000000
$ws.Range("J2").Value = 1150400
$ws.Range("K2").Value = 2000000
$ws.Range("L2").Value = 1150400
$ws.Range("M2").Value = -1999888
$ws.Range("N2").Value = -1150624
$ws.Range("H100").Value = 2134.3333
$ws.Range("I100").Value = 1200
$ws.Range("K100").Value = 2400
$ws.Range("M100").Value = -1859
$ws.Range("H117").Value = 77909
$ws.Range("J117").Value = 77909
$ws.Range("L117").Value = 77909
$ws.Range("N117").Value = -87087
$ws.Range("H132").Value = 825.3
$ws.Range("I132").Value = 802
$ws.Range("J132").Value = 868.5714
$ws.Range("K132").Value = 2406
$ws.Range("L132").Value = 2605.7142
$ws.Range("M132").Value = 124
$ws.Range("N132").Value = -7665.7142
